$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 595294.4399999999
$ws.Range("J17").Value = 595294.4399999999
$ws.Range("L17").Value = 1785883.32
$ws.Range("N17").Value = -1786219.32
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 2000
$ws.Range("K21").Value = 2000
$ws.Range("M21").Value = -1532
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("M23").Value = -1766
$ws.Range("H87").Value = 34882.645
$ws.Range("J87").Value = 34882.645
$ws.Range("L87").Value = 34882.645
$ws.Range("N87").Value = -37378.645
$ws.Range("H90").Value = 34882.645
$ws.Range("J90").Value = 34882.645
$ws.Range("L90").Value = 104647.935
$ws.Range("N90").Value = -117127.935
$ws.Range("H105").Value = 48500
$ws.Range("J105").Value = 48500
$ws.Range("L105").Value = 48500
$ws.Range("N105").Value = -55488
$ws.Range("H112").Value = 1397.159
$ws.Range("J112").Value = 1420.3489
$ws.Range("L112").Value = 4261.0467
$ws.Range("N112").Value = -6477.0467
$ws.Range("H115").Value = 4941.875
$ws.Range("I115").Value = 5259.1665
$ws.Range("J115").Value = 3990
$ws.Range("K115").Value = 15777.4995
$ws.Range("L115").Value = 11970
$ws.Range("M115").Value = -14210.4995
$ws.Range("N115").Value = -15104
$ws.Range("H132").Value = 1759.5416
$ws.Range("I132").Value = 1617.3158
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 4851.9474
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = -2321.9474
$ws.Range("N132").Value = -11960
$ws.Range("H137").Value = 3936.4412
$ws.Range("I137").Value = 1926.95
$ws.Range("K137").Value = 5780.85
$ws.Range("M137").Value = -3230.85

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6813.942
$ws.Range("I32").Value = 6271.0747
$ws.Range("K32").Value = 6271.0747
$ws.Range("M32").Value = -5984.0747
$ws.Range("H45").Value = 1489.3939
$ws.Range("I45").Value = 1516.1482
$ws.Range("J45").Value = 1369
$ws.Range("K45").Value = 1516.1482
$ws.Range("L45").Value = 1369
$ws.Range("M45").Value = -1139.1482
$ws.Range("N45").Value = -2123
$ws.Range("H61").Value = 10233.6
$ws.Range("I61").Value = 4898.2173
$ws.Range("J61").Value = 20459.75
$ws.Range("K61").Value = 4898.2173
$ws.Range("L61").Value = 20459.75
$ws.Range("M61").Value = -4686.2173
$ws.Range("N61").Value = -20883.75
$ws.Range("H74").Value = 56959.5
$ws.Range("I74").Value = 73942.8
$ws.Range("J74").Value = 8007.647
$ws.Range("K74").Value = 73942.8
$ws.Range("L74").Value = 8007.647
$ws.Range("M74").Value = -73068.8
$ws.Range("N74").Value = -9755.647000000001
$ws.Range("H77").Value = 56959.5
$ws.Range("I77").Value = 73942.8
$ws.Range("J77").Value = 8007.647
$ws.Range("K77").Value = 369714
$ws.Range("L77").Value = 40038.235
$ws.Range("M77").Value = -365346
$ws.Range("N77").Value = -48774.235
$ws.Range("H97").Value = 1054.1613
$ws.Range("I97").Value = 849.125
$ws.Range("J97").Value = 1757.1428
$ws.Range("K97").Value = 849.125
$ws.Range("L97").Value = 1757.1428
$ws.Range("M97").Value = -353.125
$ws.Range("N97").Value = -2749.1428
$ws.Range("H136").Value = 10233.6
$ws.Range("I136").Value = 4898.2173
$ws.Range("J136").Value = 20459.75
$ws.Range("K136").Value = 14694.6519
$ws.Range("L136").Value = 61379.25
$ws.Range("M136").Value = -12144.6519
$ws.Range("N136").Value = -66479.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 31826.354
$ws.Range("I134").Value = 2501.926
$ws.Range("J134").Value = 144934.86
$ws.Range("K134").Value = 7505.778
$ws.Range("L134").Value = 434804.58
$ws.Range("M134").Value = -4970.778
$ws.Range("N134").Value = -439874.58

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2130.5605
$ws.Range("I31").Value = 1505.4407
$ws.Range("J31").Value = 3283.125
$ws.Range("K31").Value = 1505.4407
$ws.Range("L31").Value = 3283.125
$ws.Range("M31").Value = -1210.4407
$ws.Range("N31").Value = -3873.125
$ws.Range("H34").Value = 2130.5605
$ws.Range("I34").Value = 1505.4407
$ws.Range("J34").Value = 3283.125
$ws.Range("K34").Value = 1505.4407
$ws.Range("L34").Value = 3283.125
$ws.Range("M34").Value = -1303.4407
$ws.Range("N34").Value = -3687.125
$ws.Range("H58").Value = 2275554.2
$ws.Range("I58").Value = 3248816.8
$ws.Range("K58").Value = 3248816.8
$ws.Range("M58").Value = -3248613.8
$ws.Range("H94").Value = 1618
$ws.Range("I94").Value = 1703
$ws.Range("J94").Value = 1504.6666
$ws.Range("K94").Value = 1703
$ws.Range("L94").Value = 1504.6666
$ws.Range("M94").Value = -1252
$ws.Range("N94").Value = -2406.6666
$ws.Range("H132").Value = 2519.4
$ws.Range("I132").Value = 1956.6086
$ws.Range("J132").Value = 3598.0833
$ws.Range("K132").Value = 5869.825800000001
$ws.Range("L132").Value = 10794.2499
$ws.Range("M132").Value = -3339.825800000001
$ws.Range("N132").Value = -15854.2499
$ws.Range("H136").Value = 2275554.2
$ws.Range("I136").Value = 3248816.8
$ws.Range("K136").Value = 9746450.399999999
$ws.Range("M136").Value = -9743900.399999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 648.06665
$ws.Range("J26").Value = 1024.7778
$ws.Range("L26").Value = 3074.3334
$ws.Range("N26").Value = -3650.3334
$ws.Range("H131").Value = 15190.123
$ws.Range("I131").Value = 390.10205
$ws.Range("J131").Value = 60515.188
$ws.Range("K131").Value = 1170.30615
$ws.Range("L131").Value = 181545.564
$ws.Range("M131").Value = 3869.69385
$ws.Range("N131").Value = -191625.564
$ws.Range("H140").Value = 1926.1154
$ws.Range("I140").Value = 1708.1364
$ws.Range("J140").Value = 3125
$ws.Range("K140").Value = 5124.4092
$ws.Range("L140").Value = 9375
$ws.Range("M140").Value = 55.59079999999994
$ws.Range("N140").Value = -19735

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 27999.5
$ws.Range("J95").Value = 27999.5
$ws.Range("L95").Value = 27999.5
$ws.Range("N95").Value = -33491.5
$ws.Range("H132").Value = 41526.895
$ws.Range("I132").Value = 85561.75
$ws.Range("J132").Value = 8500.75
$ws.Range("K132").Value = 256685.25
$ws.Range("L132").Value = 25502.25
$ws.Range("M132").Value = -254155.25
$ws.Range("N132").Value = -30562.25
$ws.Range("H135").Value = 60854
$ws.Range("J135").Value = 60854
$ws.Range("L135").Value = 60854
$ws.Range("N135").Value = -70994

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 8750
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10344
$ws.Range("H82").Value = 2599.8333
$ws.Range("I82").Value = 1411.125
$ws.Range("K82").Value = 1411.125
$ws.Range("M82").Value = -1050.125
$ws.Range("H85").Value = 2599.8333
$ws.Range("I85").Value = 1411.125
$ws.Range("K85").Value = 1411.125
$ws.Range("M85").Value = -163.125
$ws.Range("H132").Value = 3667.5625
$ws.Range("I132").Value = 3393.0833
$ws.Range("J132").Value = 4491
$ws.Range("K132").Value = 10179.2499
$ws.Range("L132").Value = 13473
$ws.Range("M132").Value = -7649.249899999999
$ws.Range("N132").Value = -18533
$ws.Range("H136").Value = 5313.9756
$ws.Range("I136").Value = 3712.8262
$ws.Range("J136").Value = 7359.8887
$ws.Range("K136").Value = 11138.4786
$ws.Range("L136").Value = 22079.6661
$ws.Range("M136").Value = -8588.4786
$ws.Range("N136").Value = -27179.6661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 43057
$ws.Range("J64").Value = 43057
$ws.Range("L64").Value = 43057
$ws.Range("N64").Value = -43553
$ws.Range("H67").Value = 43057
$ws.Range("J67").Value = 43057
$ws.Range("L67").Value = 43057
$ws.Range("N67").Value = -44773
$ws.Range("H126").Value = 1536.7273
$ws.Range("I126").Value = 1500.375
$ws.Range("J126").Value = 1633.6666
$ws.Range("K126").Value = 4501.125
$ws.Range("L126").Value = 4900.9998
$ws.Range("M126").Value = -2031.125
$ws.Range("N126").Value = -9840.9998
$ws.Range("H135").Value = 222261090
$ws.Range("J135").Value = 222261090
$ws.Range("L135").Value = 222261090
$ws.Range("N135").Value = -222271230
$ws.Range("H136").Value = 4548.9673
$ws.Range("I136").Value = 2013.0312
$ws.Range("J136").Value = 7347.241
$ws.Range("K136").Value = 6039.0936
$ws.Range("L136").Value = 22041.723
$ws.Range("M136").Value = -3489.0936
$ws.Range("N136").Value = -27141.723
